$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark that currently sits right after
#    "COS318 - FA2018" (end of the first/title paragraph).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2) Update the due-date text: "September 21th, 2017" -> "September 20th, 2018"
$d.Content.Find.Execute("September 21", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 20", 2) | Out-Null
$d.Content.Find.Execute(", 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", 2018", 2) | Out-Null

# 3) Re-create the "_GoBack" bookmark right after the new due date text
#    ("...20th, 2018"), immediately before the line break that precedes
#    "Turn in all files using Moodle".
$r = $d.Content
$r.Find.Execute(", 2018", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $insertPoint)
